# Update "Add info to Calgary (12)" - column N (CA12) values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rows that flip from 0 -> 1 in column N
$rowsToOne = @(4, 10, 19, 22, 23, 33, 36, 37, 38, 41)
foreach ($r in $rowsToOne) {
    $ws.Range("N$r").Value = 1
}

# Row that flips from 1 -> 0 in column N
$ws.Range("N24").Value = 0

# Update the active selection to match the recorded cursor position
$ws.Range("N37").Select()
